# Update cryptocurrency price (D) and volume-change (E) figures to match the
# latest scrape, preserving the original "General" cell styling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.399.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2839"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06822"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.66%  "
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.911.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07616"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.169"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "286.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.402.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007574"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.158.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9986"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.191"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.180"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.226"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.023"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1067"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05025"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7347"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9993"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02001"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.677"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.037"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8717"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.805"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "52.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +24.82%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.089"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
